# Extend the table on sheet1 with a new "2022" column (column K),
# mirroring the formatting of the existing "2021" column (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> new value for column K
$values = @{
    4  = 2022
    5  = 26.495524312074597
    6  = 59.383769502755833
    7  = 38.32334404557426
    8  = 48.136790950525594
    9  = 46.63213064070051
    10 = 32.657429481680126
    11 = 31.457245964894081
    12 = 22.734405597714229
    13 = -0.19691879995369213
    14 = 33.158040409631916
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Cells.Item($row, 10)   # column J
    $dstCell = $ws.Cells.Item($row, 11)   # column K

    # Copy formatting (number format, font, borders, fill, etc.) from J to K
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $dstCell.Value = $values[$row]
}

$excel.CutCopyMode = $false

# Update the active selection on the sheet view, as recorded in the workbook
$ws.Range("M7").Select() | Out-Null
